# "adding witch base" — populate the new "Skill" sheet with the Light/Heavy
# Magic skill-pool tables (and the class/skill-pool summary table that was
# already used on "Character"), bump Character!F13 (Wizard Shoot Power)
# 6 -> 7, and flip the active sheet/selection to match the author's saved
# view state.

$wb = $excel.ActiveWorkbook

$charSheet = $wb.Worksheets.Item("Character")
$skillSheet = $wb.Worksheets.Item("Skill")

# --- Character sheet: small stat tweak -------------------------------------
$charSheet.Range("F13").Value = 7

# --- Skill sheet: Class / Skill Pool summary (B2:F14) -----------------------
$skillSheet.Range("B2").Value = "Skill Pool"

$skillSheet.Range("B3").Value = "Class"
$skillSheet.Range("C3").Value = "Sorcerer"
$skillSheet.Range("D3").Value = "Cleric"
$skillSheet.Range("E3").Value = "Wizard"
$skillSheet.Range("F3").Value = "Druid"

$skillSheet.Range("B4").Value = "Type"
$skillSheet.Range("C4").Value = "Defender"
$skillSheet.Range("D4").Value = "Defender"
$skillSheet.Range("E4").Value = "Striker"
$skillSheet.Range("F4").Value = "Striker"

$skillSheet.Range("B6").Value = "Invisible"
$skillSheet.Range("B7").Value = "Fly"
$skillSheet.Range("B8").Value = "Heal"
$skillSheet.Range("B9").Value = "Stun"
$skillSheet.Range("B10").Value = "Damage"
$skillSheet.Range("B11").Value = "Speed"
$skillSheet.Range("B12").Value = "Tile sabotage"
$skillSheet.Range("B13").Value = "Power"
$skillSheet.Range("B14").Value = "Guard"

$skillSheet.Range("D7").Value = "B"
$skillSheet.Range("C8").Value = "B"
$skillSheet.Range("F9").Value = "B"
$skillSheet.Range("D10").Value = "C"
$skillSheet.Range("E10").Value = "B"
$skillSheet.Range("D11").Value = "C"
$skillSheet.Range("F11").Value = "B"
$skillSheet.Range("E13").Value = "B"
$skillSheet.Range("F13").Value = "B"
$skillSheet.Range("D14").Value = "B"
$skillSheet.Range("E14").Value = "B"

# "A" impact cells carry the green highlight fill (same style used on the
# existing Character!I2:M14 table).
$greenA = @("F6", "F7", "D8", "E9", "C10", "E12", "F12", "C14")
foreach ($addr in $greenA) {
    $cell = $skillSheet.Range($addr)
    $cell.Value = "A"
    $cell.Interior.Color = 5296274
}

# --- Skill sheet: Light Magic Skill Pool (I2:J8) ----------------------------
$skillSheet.Range("I2").Value = "Light Magic Skill Pool"

$skillSheet.Range("I5").Value = "Heal"
$skillSheet.Range("J5").Value = "Heal up self 3 point"

$skillSheet.Range("I6").Value = "Boost"
$skillSheet.Range("J6").Value = "Speed up self 2 point"

$skillSheet.Range("I7").Value = "Rage"
$skillSheet.Range("J7").Value = "Damage and Power up self 1 point"

$skillSheet.Range("I8").Value = "Guardian"
$skillSheet.Range("J8").Value = "Guard up self 3 point"

# --- Skill sheet: Heavy Magic Skill Pool (M2:N8) ----------------------------
$skillSheet.Range("M2").Value = "Heavy Magic Skill Pool"

$skillSheet.Range("M5").Value = "Summon Souroff"
$skillSheet.Range("N5").Value = "Summon a Souroff near the enemies"

$skillSheet.Range("M6").Value = "Holy Greave"
$skillSheet.Range("N6").Value = "Healing the team and remove stun"

$skillSheet.Range("M7").Value = "Cosmos Lighting"
$skillSheet.Range("N7").Value = "Stun all other characters"

$skillSheet.Range("M8").Value = "Nature's Command"
$skillSheet.Range("N8").Value = "Be invisible while turning all the tiles into falling rock"

# Column M needs to be wide enough to show its longest label (matches the
# author's saved best-fit width on the Heavy Magic Skill Pool column).
$skillSheet.Columns.Item(13).ColumnWidth = 20

# --- View state: Skill tab becomes the active/selected tab ------------------
$charSheet.Range("M19").Select() | Out-Null
$skillSheet.Range("I8").Select() | Out-Null
$skillSheet.Activate() | Out-Null
